# "added another year to coor ops in IC spreadsheet."
# The new year's row (r6, date 2022-01-01) is cleared back out (leaving only
# the date cell's style), and the prior year's (r5) projection values are
# updated -- net effect: CoordinatedOps becomes the active sheet/tab, with a
# new selection, and Reservoirs loses its tabSelected flag.

$wb = $excel.ActiveWorkbook

$wsReservoirs = $wb.Worksheets.Item("Reservoirs")
$wsCoord = $wb.Worksheets.Item("CoordinatedOps")

# Update CoordinatedOps row 5 (2021-01-01) projection values.
$wsCoord.Range("F5").Value = 3601
$wsCoord.Range("G5").Value = 1074

# Clear out row 6 (2022-01-01) data, leaving only A6's formatting (no value).
$wsCoord.Range("A6:G6").ClearContents()

# Reservoirs keeps its own pane/selection state, just no longer the active tab.
$wsReservoirs.Activate()
$wsReservoirs.Range("K14:K17").Select() | Out-Null

# Make CoordinatedOps the active sheet/tab, with the new selection.
$wsCoord.Activate()
$wsCoord.Range("F6").Select() | Out-Null
